# lake-ontario-spawning.xlsx - "more data cleaning and script updating"
#
# 1) The spawning-season "year" column (B) was off by one for every row -
#    bump each value up by one year.
# 2) Re-sort the data rows (A2:I13) by the date column (A) ascending - the
#    rows had fallen out of date order after the year correction.
# 3) Leave the selection where the analyst left it (E18) instead of H2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) fix the spawning-season year in column B -------------------------
for ($r = 2; $r -le 13; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value = $cell.Value2 + 1
}

# --- 2) re-sort the data by date (column A), ascending --------------------
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B1:B13")) | Out-Null
$ws.Sort.SetRange($ws.Range("A1:I13"))
$ws.Sort.Header = 1
$ws.Sort.Orientation = 1
$ws.Sort.Apply()

# --- 3) move the selection -------------------------------------------------
$ws.Range("E18").Select() | Out-Null
